$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the authoritative diff (row => column => new text value)
$updates = @(
    @{ Cell = 'D2'; Value = '58.904.66' }
    @{ Cell = 'E2'; Value = '  +2.36%  ' }
    @{ Cell = 'D3'; Value = '2.589.50' }
    @{ Cell = 'E3'; Value = '  +1.02%  ' }
    @{ Cell = 'E4'; Value = '  +0.04%  ' }
    @{ Cell = 'D5'; Value = '520.51' }
    @{ Cell = 'E5'; Value = '  +0.03%  ' }
    @{ Cell = 'D6'; Value = '139.99' }
    @{ Cell = 'E6'; Value = '  -2.60%  ' }
    @{ Cell = 'D7'; Value = '0.998' }
    @{ Cell = 'E7'; Value = '  +0.01%  ' }
    @{ Cell = 'E8'; Value = '  +0.99%  ' }
    @{ Cell = 'D9'; Value = '2.602.66' }
    @{ Cell = 'E9'; Value = '  +1.07%  ' }
    @{ Cell = 'D10'; Value = '6.55' }
    @{ Cell = 'E10'; Value = '  -0.47%  ' }
    @{ Cell = 'E11'; Value = '  +1.04%  ' }
    @{ Cell = 'E12'; Value = '  +2.17%  ' }
    @{ Cell = 'E13'; Value = '  +3.06%  ' }
    @{ Cell = 'D14'; Value = '3.047.33' }
    @{ Cell = 'E14'; Value = '  +1.07%  ' }
    @{ Cell = 'D15'; Value = '58.877.65' }
    @{ Cell = 'E15'; Value = '  +2.36%  ' }
    @{ Cell = 'D16'; Value = '20.51' }
    @{ Cell = 'E16'; Value = '  +2.01%  ' }
    @{ Cell = 'B17'; Value = 'ShibaInu' }
    @{ Cell = 'C17'; Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib' }
    @{ Cell = 'D17'; Value = '0.0000133' }
    @{ Cell = 'E17'; Value = '  +0.09%  ' }
    @{ Cell = 'B18'; Value = 'WrappedEther' }
    @{ Cell = 'C18'; Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth' }
    @{ Cell = 'D18'; Value = '2.547.57' }
    @{ Cell = 'E18'; Value = '  -0.50%  ' }
    @{ Cell = 'D19'; Value = '339.26' }
    @{ Cell = 'E19'; Value = '  +1.47%  ' }
    @{ Cell = 'D20'; Value = '4.31' }
    @{ Cell = 'E20'; Value = '  +0.82%  ' }
    @{ Cell = 'D21'; Value = '10.21' }
    @{ Cell = 'E21'; Value = '  +0.86%  ' }
    @{ Cell = 'E22'; Value = '  +4.15%  ' }
    @{ Cell = 'E23'; Value = '  +0.04%  ' }
    @{ Cell = 'D24'; Value = '66.17' }
    @{ Cell = 'E24'; Value = '  +2.47%  ' }
    @{ Cell = 'D25'; Value = '0.167' }
    @{ Cell = 'E25'; Value = '  +0.55%  ' }
    @{ Cell = 'D26'; Value = '0.405' }
    @{ Cell = 'E26'; Value = '  +1.34%  ' }
    @{ Cell = 'E27'; Value = '  +0.13%  ' }
    @{ Cell = 'D28'; Value = '7.06' }
    @{ Cell = 'E28'; Value = '  +1.68%  ' }
    @{ Cell = 'D29'; Value = '0.998' }
    @{ Cell = 'E29'; Value = '  +0.04%  ' }
    @{ Cell = 'E30'; Value = '  -2.37%  ' }
    @{ Cell = 'D31'; Value = '5.94' }
    @{ Cell = 'E31'; Value = '  -4.86%  ' }
    @{ Cell = 'E32'; Value = '  -0.92%  ' }
    @{ Cell = 'E33'; Value = '  +1.29%  ' }
    @{ Cell = 'D34'; Value = '148.84' }
    @{ Cell = 'E34'; Value = '  +0.37%  ' }
    @{ Cell = 'D35'; Value = '3.98' }
    @{ Cell = 'E35'; Value = '  -0.78%  ' }
    @{ Cell = 'E36'; Value = '  -0.99%  ' }
    @{ Cell = 'B37'; Value = 'Stacks' }
    @{ Cell = 'C37'; Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx' }
    @{ Cell = 'D37'; Value = '1.47' }
    @{ Cell = 'E37'; Value = '  +1.96%  ' }
    @{ Cell = 'B38'; Value = 'OKB' }
    @{ Cell = 'C38'; Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb' }
    @{ Cell = 'D38'; Value = '36.27' }
    @{ Cell = 'E38'; Value = '  +1.14%  ' }
    @{ Cell = 'E39'; Value = '  -0.12%  ' }
    @{ Cell = 'E40'; Value = '  -1.92%  ' }
    @{ Cell = 'E41'; Value = '  +0.49%  ' }
    @{ Cell = 'D43'; Value = '274.88' }
    @{ Cell = 'E43'; Value = '  +2.31%  ' }
    @{ Cell = 'E44'; Value = '  +1.05%  ' }
    @{ Cell = 'D45'; Value = '0.590' }
    @{ Cell = 'E45'; Value = '  +0.57%  ' }
    @{ Cell = 'D46'; Value = '0.0952' }
    @{ Cell = 'E46'; Value = '  -0.07%  ' }
    @{ Cell = 'D47'; Value = '0.0522' }
    @{ Cell = 'E47'; Value = '  +0.54%  ' }
    @{ Cell = 'E48'; Value = '  -0.99%  ' }
    @{ Cell = 'D49'; Value = '1.992.89' }
    @{ Cell = 'E49'; Value = '  +1.05%  ' }
    @{ Cell = 'B50'; Value = 'VeChain' }
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet' }
    @{ Cell = 'D50'; Value = '0.0221' }
    @{ Cell = 'E50'; Value = '  +0.88%  ' }
    @{ Cell = 'B51'; Value = 'RenderToken' }
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr' }
    @{ Cell = 'D51'; Value = '4.50' }
    @{ Cell = 'E51'; Value = '  -4.68%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}
